# Updated all assays to accommodate the newly introduced dataset type
#
# 1. dataset_type list: remove "nanoPOTS" and "NanoDESI", add a new
#    "2D Imaging Mass Cytometry" entry right after "MALDI".
# 2. acquisition_instrument_model list: add "STELLARIS 5" right after
#    "SCN400" and "Unknown" right after
#    "Resolve Biosciences Molecular Cartography".
# 3. .metadata sheet: bump the pav:createdOn timestamp.
# 4. Repair the two data-validation list ranges on the MALDI sheet so
#    they keep pointing at the full (resized) lookup ranges.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. dataset_type
# ---------------------------------------------------------------------
$dsType = $wb.Worksheets.Item("dataset_type")

$nanoPOTS = $dsType.Columns.Item(1).Find("nanoPOTS")
$dsType.Rows.Item($nanoPOTS.Row).Delete()

$nanoDESI = $dsType.Columns.Item(1).Find("NanoDESI")
$dsType.Rows.Item($nanoDESI.Row).Delete()

$maldiRow = $dsType.Columns.Item(1).Find("MALDI")
$newRow = $maldiRow.Row + 1
$dsType.Rows.Item($newRow).Insert()
$dsType.Cells.Item($newRow, 1).Value = "2D Imaging Mass Cytometry"
$dsType.Cells.Item($newRow, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296"

$dsTypeCount = $dsType.UsedRange.Rows.Count

# ---------------------------------------------------------------------
# 2. acquisition_instrument_model
# ---------------------------------------------------------------------
$acqModel = $wb.Worksheets.Item("acquisition_instrument_model")

$scn400 = $acqModel.Columns.Item(1).Find("SCN400")
$stellarisRow = $scn400.Row + 1
$acqModel.Rows.Item($stellarisRow).Insert()
$acqModel.Cells.Item($stellarisRow, 1).Value = "STELLARIS 5"
$acqModel.Cells.Item($stellarisRow, 2).Value = "https://identifiers.org/RRID:SCR_024663"

$resolveMC = $acqModel.Columns.Item(1).Find("Resolve Biosciences Molecular Cartography")
$unknownRow = $resolveMC.Row + 1
$acqModel.Rows.Item($unknownRow).Insert()
$acqModel.Cells.Item($unknownRow, 1).Value = "Unknown"
$acqModel.Cells.Item($unknownRow, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998"

$acqModelCount = $acqModel.UsedRange.Rows.Count

# ---------------------------------------------------------------------
# 3. .metadata : pav:createdOn
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item(".metadata")
$createdOnHeader = $meta.Rows.Item(1).Find("pav:createdOn")
$createdOnCol = $createdOnHeader.Column
$meta.Cells.Item(2, $createdOnCol).Value = "2023-11-02T15:46:29-07:00"

# ---------------------------------------------------------------------
# 4. Fix up the data validation ranges on the main MALDI sheet
# ---------------------------------------------------------------------
$main = $wb.Worksheets.Item("MALDI")

$dsRange = $main.Range("D2:D1001")
$dsValidation = $dsRange.Validation
$dsValidation.Delete()
$dsRange.Validation.Add(3, 1, 1, "='dataset_type'!`$A`$1:`$A`$$dsTypeCount")
$dsValidation = $dsRange.Validation
$dsValidation.ErrorTitle = "Validation Error"
$dsValidation.ErrorMessage = ""
$dsValidation.IgnoreBlank = $true
$dsValidation.ShowError = $true
$dsValidation.ShowInput = $false

$modelRange = $main.Range("H2:H1001")
$modelValidation = $modelRange.Validation
$modelValidation.Delete()
$modelRange.Validation.Add(3, 1, 1, "='acquisition_instrument_model'!`$A`$1:`$A`$$acqModelCount")
$modelValidation = $modelRange.Validation
$modelValidation.ErrorTitle = "Validation Error"
$modelValidation.ErrorMessage = ""
$modelValidation.IgnoreBlank = $true
$modelValidation.ShowError = $true
$modelValidation.ShowInput = $false

Write-Output "dataset_type rows: $dsTypeCount"
Write-Output "acquisition_instrument_model rows: $acqModelCount"
